$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.7103813333333333
$ws.Range("H2").Value = 2.131144
$ws.Range("I2").Value = 0.7576743564291667
$ws.Range("J2").Value = 0.7576743564291667
$ws.Range("K2").Value = 2.0
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3659943333333333
$ws.Range("N2").Value = 1.097983
$ws.Range("O2").Value = 0.006726051721149161
$ws.Range("P2").Value = 0.006726051721149162
$ws.Range("Q2").Value = 0.2599955425057778
$ws.Range("R2").Value = 2.339959882552
$ws.Range("S2").Value = 0.00509615690913098
$ws.Range("T2").Value = 0.005096156909130981

$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.7103813333333333
$ws.Range("H3").Value = 2.131144
$ws.Range("I3").Value = 0.7576743564291667
$ws.Range("J3").Value = 0.7576743564291667
$ws.Range("K3").Value = 1.0
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09131133333333334
$ws.Range("N3").Value = 0.273934
$ws.Range("O3").Value = 0.001678071748088335
$ws.Range("P3").Value = 0.001678071748088335
$ws.Range("Q3").Value = 0.06486586672177778
$ws.Range("R3").Value = 0.583792800496
$ws.Range("S3").Value = 0.001271431931774796
$ws.Range("T3").Value = 0.001271431931774796

$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.7103813333333333
$ws.Range("H4").Value = 2.131144
$ws.Range("I4").Value = 0.7576743564291667
$ws.Range("J4").Value = 0.7576743564291667
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 53.897087
$ws.Range("N4").Value = 161.691261
$ws.Range("O4").Value = 0.9904923704135933
$ws.Range("P4").Value = 0.9904923704135934
$ws.Range("Q4").Value = 38.28748452584266
$ws.Range("R4").Value = 344.587360732584
$ws.Range("S4").Value = 0.750470669301119
$ws.Range("T4").Value = 0.7504706693011192

$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.7103813333333333
$ws.Range("H5").Value = 2.131144
$ws.Range("I5").Value = 0.7576743564291667
$ws.Range("J5").Value = 0.7576743564291667
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06004666666666667
$ws.Range("N5").Value = 0.18014
$ws.Range("O5").Value = 0.001103506117169219
$ws.Range("P5").Value = 0.001103506117169219
$ws.Range("Q5").Value = 0.04265603112888889
$ws.Range("R5").Value = 0.3839042801600001
$ws.Range("S5").Value = 0.0008360982871418363
$ws.Range("T5").Value = 0.0008360982871418364

$ws.Range("E6").Value = 2.0
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2272
$ws.Range("H6").Value = 0.6816
$ws.Range("I6").Value = 0.2423256435708333
$ws.Range("J6").Value = 0.2423256435708333
$ws.Range("K6").Value = 2.0
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3659943333333333
$ws.Range("N6").Value = 1.097983
$ws.Range("O6").Value = 0.006726051721149161
$ws.Range("P6").Value = 0.006726051721149162
$ws.Range("Q6").Value = 0.08315391253333332
$ws.Range("R6").Value = 0.7483852127999999
$ws.Range("S6").Value = 0.001629894812018182
$ws.Range("T6").Value = 0.001629894812018182

$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2272
$ws.Range("H7").Value = 0.6816
$ws.Range("I7").Value = 0.2423256435708333
$ws.Range("J7").Value = 0.2423256435708333
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.09131133333333334
$ws.Range("N7").Value = 0.273934
$ws.Range("O7").Value = 0.001678071748088335
$ws.Range("P7").Value = 0.001678071748088335
$ws.Range("Q7").Value = 0.02074593493333333
$ws.Range("R7").Value = 0.1867134144
$ws.Range("S7").Value = 0.0004066398163135392
$ws.Range("T7").Value = 0.0004066398163135392

$ws.Range("E8").Value = 2.0
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2272
$ws.Range("H8").Value = 0.6816
$ws.Range("I8").Value = 0.2423256435708333
$ws.Range("J8").Value = 0.2423256435708333
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 53.897087
$ws.Range("N8").Value = 161.691261
$ws.Range("O8").Value = 0.9904923704135933
$ws.Range("P8").Value = 0.9904923704135934
$ws.Range("Q8").Value = 12.2454181664
$ws.Range("R8").Value = 110.2087634976
$ws.Range("S8").Value = 0.2400217011124742
$ws.Range("T8").Value = 0.2400217011124743

$ws.Range("E9").Value = 2.0
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2272
$ws.Range("H9").Value = 0.6816
$ws.Range("I9").Value = 0.2423256435708333
$ws.Range("J9").Value = 0.2423256435708333
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06004666666666667
$ws.Range("N9").Value = 0.18014
$ws.Range("O9").Value = 0.001103506117169219
$ws.Range("P9").Value = 0.001103506117169219
$ws.Range("Q9").Value = 0.01364260266666667
$ws.Range("R9").Value = 0.122783424
$ws.Range("S9").Value = 0.0002674078300273823
$ws.Range("T9").Value = 0.0002674078300273823
